# Add rank-of-AICc column (H) per block of 10 replicate rows, and a small
# summary table (J2:K5) that sums the rank-1 counts for a few selected
# models (gaussian_1987, lrf_1991, modifiedgaussian_2006, pawar_2018).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H column: RANK(F<row>, F<blockStart>:F<blockEnd>, 1) for each block ---
$groupStarts = 2,12,22,32,42,52,62,72,82
foreach ($s in $groupStarts) {
    $e = $s + 9
    for ($r = $s; $r -le $e; $r++) {
        $ws.Range("H$r").Formula = "=RANK(F$r,F$($s):F$($e),1)"
    }
}

# New H cells should not inherit the coloured "customFormat" row styling -
# reset them back to the plain/default style.
$ws.Range("H2:H91").Style = "Normal"

# --- J2:J5 / K2:K5 summary table ---
$ws.Range("J2").Value = "gaussian_1987"
$ws.Range("J3").Value = "lrf_1991"
$ws.Range("J4").Value = "modifiedgaussian_2006"
$ws.Range("J5").Value = "pawar_2018"

# Match the shaded-row styling used elsewhere in the sheet: J2 picks up the
# accent used on row 86, J3:J5 pick up the accent used on row 2.
$ws.Range("A86").Copy()
$ws.Range("J2").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("J3:J5").PasteSpecial(-4122)

$ws.Range("K2").Formula = "=SUM(H6,H16,H26,H36,H46,H56,H66,H76,H86)"
$ws.Range("K3").Formula = "=SUM(H9,H19,H29,H39,H49,H59,H69,H79,H89)"
$ws.Range("K4").Formula = "=SUM(H2,H12,H22,H32,H42,H52,H62,H72,H82)"
$ws.Range("K5").Formula = "=SUM(H10,H20,H30,H40,H50,H60,H70,H80,H90)"

$ws.Range("J2").Select()
